# "Fruta / hortaliza, semanal" - add a new weekly price record for
# Agrícola del Norte S.A. de Arica / Zapallo / Camote.
#
# The new record is inserted as row 26 (pushing the existing rows 26-40
# down to 27-41), extending the sheet's used range from A1:R40 to A1:R41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 26 - this shifts rows
# 26-40 down to 27-41 and extends the sheet dimension accordingly.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44873
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112045
$ws.Cells.Item(26, 7).Value = "Zapallo"
$ws.Cells.Item(26, 8).Value = "Camote"
$ws.Cells.Item(26, 9).Value = "1a nueva(o)"
$ws.Cells.Item(26, 10).Value = 800
$ws.Cells.Item(26, 11).Value = 950
$ws.Cells.Item(26, 12).Value = 980
$ws.Cells.Item(26, 13).Value = 965
$ws.Cells.Item(26, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 965
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
